$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 6) with Item ID and Item Name, matching existing data pattern.
# Force column A to be stored as text so "2830113" is kept as a shared string, not a number.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2830113"
$ws.Range("B6").Value = "Le Chat HS 300 gr"

# Match formatting of the row above (border style) for the new row A6:D6
$ws.Range("A5:D5").Copy() | Out-Null
$ws.Range("A6:D6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the values since paste-special(formats) should not touch them, but ensure correctness
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2830113"
$ws.Range("B6").Value = "Le Chat HS 300 gr"

# Update the active selection to reflect new last cell, similar to the diff (A6 selected)
$ws.Range("A6").Select() | Out-Null

$wb.Save()
